# Update "Ireland-manual" sheet with data for 2020-04-17 .. 2020-04-26
# (commit message: "Update to data for April 26th")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ireland-manual")
$ws.Activate()

# Each entry: row, CountryRegion, Lat, Long, date, cases, type
$rows = @(
    @(174, "Ireland", 53.1424, -7.6921, "2020-04-17", 597, "confirmed"),
    @(175, "Ireland", 53.1424, -7.6921, "2020-04-17", 44,  "death"),
    @(176, "Ireland", 53.1424, -7.6921, "2020-04-18", 630, "confirmed"),
    @(177, "Ireland", 53.1424, -7.6921, "2020-04-18", 41,  "death"),
    @(178, "Ireland", 53.1424, -7.6921, "2020-04-19", 445, "confirmed"),
    @(179, "Ireland", 53.1424, -7.6921, "2020-04-19", 39,  "death"),
    @(180, "Ireland", 53.1424, -7.6921, "2020-04-20", 401, "confirmed"),
    @(181, "Ireland", 53.1424, -7.6921, "2020-04-20", 77,  "death"),
    @(182, "Ireland", 53.1424, -7.6921, "2020-04-21", 388, "confirmed"),
    @(183, "Ireland", 53.1424, -7.6921, "2020-04-21", 44,  "death"),
    @(184, "Ireland", 53.1424, -7.6921, "2020-04-22", 631, "confirmed"),
    @(185, "Ireland", 53.1424, -7.6921, "2020-04-22", 49,  "death"),
    @(186, "Ireland", 53.1424, -7.6921, "2020-04-23", 936, "confirmed"),
    @(187, "Ireland", 53.1424, -7.6921, "2020-04-23", 28,  "death"),
    @(188, "Ireland", 53.1424, -7.6921, "2020-04-24", 577, "confirmed"),
    @(189, "Ireland", 53.1424, -7.6921, "2020-04-24", 37,  "death"),
    @(190, "Ireland", 53.1424, -7.6921, "2020-04-25", 377, "confirmed"),
    @(191, "Ireland", 53.1424, -7.6921, "2020-04-25", 42,  "death"),
    @(192, "Ireland", 53.1424, -7.6921, "2020-04-26", 701, "confirmed"),
    @(193, "Ireland", 53.1424, -7.6921, "2020-04-26", 26,  "death")
)

foreach ($r in $rows) {
    $rowNum = $r[0]

    $ws.Cells.Item($rowNum, 2).Value = $r[1]   # B - Country.Region
    $ws.Cells.Item($rowNum, 3).Value = $r[2]   # C - Lat
    $ws.Cells.Item($rowNum, 4).Value = $r[3]   # D - Long

    $dateCell = $ws.Cells.Item($rowNum, 5)     # E - date (stored as text)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $r[4]

    $ws.Cells.Item($rowNum, 6).Value = $r[5]   # F - cases
    $ws.Cells.Item($rowNum, 7).Value = $r[6]   # G - type
}

# Match the author's final selection / scroll position on the sheet
$ws.Range("F192").Select() | Out-Null
